$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column F ("jezyk" / language column). Columns G:J shift
# left to F:I, the now-unused shared string is dropped, and the sheet
# dimension shrinks from A1:J3 to A1:I3.
$ws.Range("F1").EntireColumn.Delete()

# Preserve the current (correct) formatting of the cells that now hold the
# e-mail addresses (H2:H3, former I2:I3 which used the blue hyperlink font)
# before we touch the Hyperlinks collection.
$mailColor = $ws.Range("H2:H3").Font.Color
$mailName = $ws.Range("H2:H3").Font.Name
$mailSize = $ws.Range("H2:H3").Font.Size

# The engine does not re-target Hyperlinks automatically when columns shift,
# so the old hyperlinks still point at I2/I3. Drop them and recreate them on
# the cells that now actually hold the e-mail addresses (H2/H3).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:olo@gmail.com", "", "", "olo@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:ala@gmail.com", "", "", "ala@gmail.com")

# Restore the original cell formatting (Hyperlinks.Add resets it to the
# default blue/underlined "Hyperlink" style).
$ws.Range("H2:H3").Font.Color = $mailColor
$ws.Range("H2:H3").Font.Name = $mailName
$ws.Range("H2:H3").Font.Size = $mailSize
$ws.Range("H2:H3").Font.Underline = -4142

# Match the new active selection left behind by deleting column F.
[void]$ws.Range("F1").Select()
